$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column I (SKLearn Predictions US) and J (Statsmodel Predictions US)
$updates = @(
    ,@(2, 89680.84809134352, 43959.32786496047)
    ,@(3, 627538.602136745, 607611.0254131673)
    ,@(4, 254093.6247243253, 250173.3809718147)
    ,@(5, 3414462.211413235, 3431326.179150726)
    ,@(6, 501727.6568343503, 476801.1849629739)
    ,@(7, 317376.9762800836, 340215.5535873612)
    ,@(8, 99263.76315908029, 68066.1379514771)
    ,@(9, 1815240.234505753, 1824997.720047017)
    ,@(10, 864662.272023945, 871773.2534000638)
    ,@(11, 134842.4186845938, 92570.17244002626)
    ,@(12, 165115.1453352451, 116866.0726362292)
    ,@(13, 1120105.251896431, 1120985.437135335)
    ,@(14, 590150.8272218243, 568202.3572237636)
    ,@(15, 279281.9469904518, 240014.9141841478)
    ,@(16, 255440.6568050117, 224667.3468455353)
    ,@(17, 373879.7596804706, 371826.6473182739)
    ,@(18, 396597.9496487469, 407462.4750298752)
    ,@(19, 130003.7145819376, 96932.0789611509)
    ,@(20, 538054.014029504, 513238.5174013065)
    ,@(21, 668638.8956284039, 663963.148942236)
    ,@(22, 854467.8260094296, 839963.9496432851)
    ,@(23, 504615.4309439024, 466074.092693327)
    ,@(24, 234478.6499564599, 247490.7002488937)
    ,@(25, 510394.7520563612, 495937.603385362)
    ,@(26, 99462.33297687286, 78315.38847442208)
    ,@(27, 185394.5360685421, 142752.7195299924)
    ,@(28, 260812.5676274506, 252643.6218806413)
    ,@(29, 133295.0358113525, 91760.63520170399)
    ,@(30, 767158.0883603166, 767238.5365156317)
    ,@(31, 184138.7927124809, 181235.7177124913)
    ,@(32, 1800762.444166858, 1831426.851220681)
    ,@(33, 169114.8240730867, 156316.7480038513)
    ,@(34, 81170.8125699799, 55090.23943158255)
    ,@(35, 977968.8954195726, 967994.62904577)
    ,@(36, 332496.2263751038, 328006.3315225848)
    ,@(37, 365172.7575407652, 334976.5029625194)
    ,@(38, 1061039.27273517, 1062841.283519712)
    ,@(39, 212459.9106142865, 292323.6536746015)
    ,@(40, 115355.960328244, 95874.65483752887)
    ,@(41, 435260.106267967, 432117.6860542103)
    ,@(42, 92741.40835667176, 45861.73425586108)
    ,@(43, 578047.5756680272, 575525.6843469266)
    ,@(44, 2362087.553650852, 2376344.81740964)
    ,@(45, 291179.9493046599, 235718.6847034938)
    ,@(46, 71854.70176299234, 36179.27193524612)
    ,@(47, 701168.1783848826, 694942.5361249596)
    ,@(48, 630938.0565473718, 612272.0202586612)
    ,@(49, 162206.7565360145, 145365.2254907412)
    ,@(50, 520004.7755928079, 477394.3779992069)
    ,@(51, 69731.05191004163, 21443.52047887332)
)

foreach ($u in $updates) {
    $row = $u[0]
    $ws.Cells.Item($row, 9).Value = $u[1]
    $ws.Cells.Item($row, 10).Value = $u[2]
}